# Remove all slide-level animation timing (the <p:timing> tree) from the
# slides that have it. This strips every effect out of each slide's main
# animation sequence, which drops the now-empty <p:timing> element from
# the underlying XML altogether.

$p = $ppt.ActivePresentation

for ($slideIdx = 1; $slideIdx -le $p.Slides.Count; $slideIdx++) {
    $s = $p.Slides.Item($slideIdx)
    $seq = $s.TimeLine.MainSequence
    for ($i = $seq.Count; $i -ge 1; $i--) {
        $seq.Item($i).Delete()
    }
}
